{"js": "// ENGR Prefix in License Signatory\n// Update the addressee block, salutation, body amounts and the\n// \"settle the amount of\" figure to reflect the new signatory\n// (Davao Sugar Central Company, Inc.) and new fee figures.\n\nconst body = context.document.body;\n\n// Simple, short, unique strings can be handled with body.search() +\n// Range.insertText(..., \"Replace\") which is a format-preserving\n// text replace (keeps the run's bold/font/size properties intact).\nconst replacements = [\n  [\"MS. MINNIE O. CHUA\", \"MR. JONATHAN T. GOTIANUN\"],\n  [\"President & COO\", \"President\"],\n  [\"Victorias Milling Company, Inc.\", \"Davao Sugar Central Company, Inc.\"],\n  [\"Ossorio St., Brgy. XVI, Victorias City, Neg. Occ.\", \"5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City\"],\n  [\"Dear Pres. Chua:\", \"Dear Mr. Gotianun:\"],\n  [\"THIRTEEN THOUSAND FIVE HUNDRED FIFTY  PESOS (PHP 13,550.00)\", \"EIGHT HUNDRED  PESOS (PHP 800.00)\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const found = body.search(find, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n  for (const rng of found.items) {\n    rng.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n\n// The production-estimate paragraph is a single long run (> 255 chars),\n// which exceeds Word's search-string limit, so address it directly via\n// its run text instead of body.search().\nconst bigFind =\n  \"Please be informed that based on your submitted production estimate of 263,250.00 Metric Tons or 5,265,000.00 Lkg., your Milling License Fee for Crop Year 2020 - 2021 is FOURTEEN THOUSAND FIVE HUNDRED  (PHP 14,500.00) PESOS.  However, you have an excess payment in your Milling License Fee for CY 2020 - 2021 in the amount of NINE HUNDRED FIFTY  PESOS (PHP 950.00).\";\nconst bigReplace =\n  \"Please be informed that based on your submitted production estimate of 1,000.00 Metric Tons or 1,000.00 Lkg., your Milling License Fee for Crop Year 2020 - 2021 is ONE THOUSAND  (PHP 1,000.00) PESOS.  However, you have an excess payment in your Milling License Fee for CY 2020 - 2021 in the amount of TWO HUNDRED  PESOS (PHP 200.00).\";\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet bigPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Please be informed that based on your submitted production estimate\") === 0) {\n    bigPara = p;\n    break;\n  }\n}\nif (!bigPara) {\n  throw new Error(\"Production-estimate paragraph not found\");\n}\nbigPara.insertText(bigReplace, \"Replace\");\nawait context.sync();\n", "ps1": "# ENGR Prefix in License Signatory\n# Update the addressee block, salutation, body amounts and the\n# \"settle the amount of\" figure to reflect the new signatory\n# (Davao Sugar Central Company, Inc.) and new fee figures.\n\n$d = $word.ActiveDocument\n\n# Short, unique strings -> Find/Replace across the whole document body.\n$replacements = @(\n    @(\"MS. MINNIE O. CHUA\", \"MR. JONATHAN T. GOTIANUN\"),\n    @(\"President & COO\", \"President\"),\n    @(\"Victorias Milling Company, Inc.\", \"Davao Sugar Central Company, Inc.\"),\n    @(\"Ossorio St., Brgy. XVI, Victorias City, Neg. Occ.\", \"5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City\"),\n    @(\"Dear Pres. Chua:\", \"Dear Mr. Gotianun:\"),\n    @(\"THIRTEEN THOUSAND FIVE HUNDRED FIFTY  PESOS (PHP 13,550.00)\", \"EIGHT HUNDRED  PESOS (PHP 800.00)\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $findText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $replaceText,\n        2\n    )\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# The production-estimate paragraph is a single long run (> 255 chars),\n# which exceeds Word's Find-text limit, so address it directly via the\n# paragraph's Range instead of Find.Execute.\n$bigReplace = \"Please be informed that based on your submitted production estimate of 1,000.00 Metric Tons or 1,000.00 Lkg., your Milling License Fee for Crop Year 2020 - 2021 is ONE THOUSAND  (PHP 1,000.00) PESOS.  However, you have an excess payment in your Milling License Fee for CY 2020 - 2021 in the amount of TWO HUNDRED  PESOS (PHP 200.00).\"\n\n$bigParaFound = $false\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"Please be informed that based on your submitted production estimate\")) {\n        $r = $p.Range\n        $r.End = $r.End - 1\n        $r.Text = $bigReplace\n        $bigParaFound = $true\n        break\n    }\n}\nif (-not $bigParaFound) {\n    throw \"Production-estimate paragraph not found\"\n}\n"}
